$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.201571333333333
$ws.Range("H2").Value = 18.604714
$ws.Range("I2").Value = 0.05221490529364391
$ws.Range("J2").Value = 0.07406232529850043
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.211442
$ws.Range("N2").Value = 99.634326
$ws.Range("O2").Value = 0.211580186305583
$ws.Range("P2").Value = 0.2175281749633597
$ws.Range("Q2").Value = 205.9631266458626
$ws.Range("R2").Value = 1853.668139812764
$ws.Range("S2").Value = 0.01104763938995755
$ws.Range("T2").Value = 0.01611064245572546

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.201571333333333
$ws.Range("H3").Value = 18.604714
$ws.Range("I3").Value = 0.05221490529364391
$ws.Range("J3").Value = 0.07406232529850043
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 51.17424933333334
$ws.Range("N3").Value = 153.522748
$ws.Range("O3").Value = 0.3260158715178649
$ws.Range("P3").Value = 0.3351809012869699
$ws.Range("Q3").Value = 317.3607576704524
$ws.Range("R3").Value = 2856.246819034072
$ws.Range("S3").Value = 0.0170228878555301
$ws.Range("T3").Value = 0.02482427694496013

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.201571333333333
$ws.Range("H4").Value = 18.604714
$ws.Range("I4").Value = 0.05221490529364391
$ws.Range("J4").Value = 0.07406232529850043
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 22.19450366666666
$ws.Range("N4").Value = 66.58351099999999
$ws.Range("O4").Value = 0.1413945597650736
$ws.Range("P4").Value = 0.1453694746776606
$ws.Range("Q4").Value = 137.6407976967615
$ws.Range("R4").Value = 1238.767179270854
$ws.Range("S4").Value = 0.007382903547169794
$ws.Range("T4").Value = 0.01076640132204902

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.201571333333333
$ws.Range("H5").Value = 18.604714
$ws.Range("I5").Value = 0.05221490529364391
$ws.Range("J5").Value = 0.07406232529850043
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.51216133333333
$ws.Range("N5").Value = 112.536484
$ws.Range("O5").Value = 0.2389787857941174
$ws.Range("P5").Value = 0.2456970098971044
$ws.Range("Q5").Value = 232.6343443761751
$ws.Range("R5").Value = 2093.709099385576
$ws.Range("S5").Value = 0.01247825466742985
$ws.Range("T5").Value = 0.01819689187186823

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.201571333333333
$ws.Range("H6").Value = 18.604714
$ws.Range("I6").Value = 0.05221490529364391
$ws.Range("J6").Value = 0.07406232529850043
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 12.8762265
$ws.Range("N6").Value = 25.752453
$ws.Range("O6").Value = 0.08203059661736112
$ws.Range("P6").Value = 0.05622443917490542
$ws.Range("Q6").Value = 79.852837143907
$ws.Range("R6").Value = 479.117022863442
$ws.Range("S6").Value = 0.004283219833556617
$ws.Range("T6").Value = 0.004164112703897596

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.461641333333333
$ws.Range("H7").Value = 22.384924
$ws.Range("I7").Value = 0.06282422221945559
$ws.Range("J7").Value = 0.0891107233935555
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.211442
$ws.Range("N7").Value = 99.634326
$ws.Range("O7").Value = 0.211580186305583
$ws.Range("P7").Value = 0.2175281749633597
$ws.Range("Q7").Value = 247.8118683668026
$ws.Range("R7").Value = 2230.306815301224
$ws.Range("S7").Value = 0.01329236064169576
$ws.Range("T7").Value = 0.01938409302946489

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.461641333333333
$ws.Range("H8").Value = 22.384924
$ws.Range("I8").Value = 0.06282422221945559
$ws.Range("J8").Value = 0.0891107233935555
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 51.17424933333334
$ws.Range("N8").Value = 153.522748
$ws.Range("O8").Value = 0.3260158715178649
$ws.Range("P8").Value = 0.3351809012869699
$ws.Range("Q8").Value = 381.8438940279058
$ws.Range("R8").Value = 3436.595046251152
$ws.Range("S8").Value = 0.02048169355930783
$ws.Range("T8").Value = 0.02986821258138581

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.461641333333333
$ws.Range("H9").Value = 22.384924
$ws.Range("I9").Value = 0.06282422221945559
$ws.Range("J9").Value = 0.0891107233935555
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.19450366666666
$ws.Range("N9").Value = 66.58351099999999
$ws.Range("O9").Value = 0.1413945597650736
$ws.Range("P9").Value = 0.1453694746776606
$ws.Range("Q9").Value = 165.6074259320182
$ws.Range("R9").Value = 1490.466833388164
$ws.Range("S9").Value = 0.008883003243303081
$ws.Range("T9").Value = 0.01295397904786748

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.461641333333333
$ws.Range("H10").Value = 22.384924
$ws.Range("I10").Value = 0.06282422221945559
$ws.Range("J10").Value = 0.0891107233935555
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 37.51216133333333
$ws.Range("N10").Value = 112.536484
$ws.Range("O10").Value = 0.2389787857941174
$ws.Range("P10").Value = 0.2456970098971044
$ws.Range("Q10").Value = 279.9022935074684
$ws.Range("R10").Value = 2519.120641567216
$ws.Range("S10").Value = 0.01501365634446531
$ws.Range("T10").Value = 0.02189423828756454

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.461641333333333
$ws.Range("H11").Value = 22.384924
$ws.Range("I11").Value = 0.06282422221945559
$ws.Range("J11").Value = 0.0891107233935555
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 12.8762265
$ws.Range("N11").Value = 25.752453
$ws.Range("O11").Value = 0.08203059661736112
$ws.Range("P11").Value = 0.05622443917490542
$ws.Range("Q11").Value = 96.077783869762
$ws.Range("R11").Value = 576.466703218572
$ws.Range("S11").Value = 0.005153508430683617
$ws.Range("T11").Value = 0.005010200447272782

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 105.106922
$ws.Range("H12").Value = 210.213844
$ws.Range("I12").Value = 0.8849608724869005
$ws.Range("J12").Value = 0.836826951307944
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 33.211442
$ws.Range("N12").Value = 99.634326
$ws.Range("O12").Value = 0.211580186305583
$ws.Range("P12").Value = 0.2175281749633597
$ws.Range("Q12").Value = 3490.752443801524
$ws.Range("R12").Value = 20944.51466280914
$ws.Range("S12").Value = 0.1872401862739297
$ws.Range("T12").Value = 0.1820334394781693

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 105.106922
$ws.Range("H13").Value = 210.213844
$ws.Range("I13").Value = 0.8849608724869005
$ws.Range("J13").Value = 0.836826951307944
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 51.17424933333334
$ws.Range("N13").Value = 153.522748
$ws.Range("O13").Value = 0.3260158715178649
$ws.Range("P13").Value = 0.3351809012869699
$ws.Range("Q13").Value = 5378.767833087219
$ws.Range("R13").Value = 32272.60699852331
$ws.Range("S13").Value = 0.288511290103027
$ws.Range("T13").Value = 0.280488411760624

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 105.106922
$ws.Range("H14").Value = 210.213844
$ws.Range("I14").Value = 0.8849608724869005
$ws.Range("J14").Value = 0.836826951307944
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 22.19450366666666
$ws.Range("N14").Value = 66.58351099999999
$ws.Range("O14").Value = 0.1413945597650736
$ws.Range("P14").Value = 0.1453694746776606
$ws.Range("Q14").Value = 2332.795965721047
$ws.Range("R14").Value = 13996.77579432628
$ws.Range("S14").Value = 0.1251286529746008
$ws.Range("T14").Value = 0.1216490943077441

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 105.106922
$ws.Range("H15").Value = 210.213844
$ws.Range("I15").Value = 0.8849608724869005
$ws.Range("J15").Value = 0.836826951307944
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 37.51216133333333
$ws.Range("N15").Value = 112.536484
$ws.Range("O15").Value = 0.2389787857941174
$ws.Range("P15").Value = 0.2456970098971044
$ws.Range("Q15").Value = 3942.787815314082
$ws.Range("R15").Value = 23656.7268918845
$ws.Range("S15").Value = 0.2114868747822222
$ws.Range("T15").Value = 0.2056058797376716

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 105.106922
$ws.Range("H16").Value = 210.213844
$ws.Range("I16").Value = 0.8849608724869005
$ws.Range("J16").Value = 0.836826951307944
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 12.8762265
$ws.Range("N16").Value = 25.752453
$ws.Range("O16").Value = 0.08203059661736112
$ws.Range("P16").Value = 0.05622443917490542
$ws.Range("Q16").Value = 1353.380534389833
$ws.Range("R16").Value = 5413.522137559333
$ws.Range("S16").Value = 0.07259386835312089
$ws.Range("T16").Value = 0.04705012602373504
